$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Lake Inez" row at the bottom of the data table (row 20)
$ws.Range("A20").Value = "Lake Inez"
$ws.Range("B20").Value = "Leslieville"
$ws.Range("C20").Value = "Secret patio menu, really anything on the menu is going to be outrageous; as of August 2024, our favorite restaurant in Toronto"
$ws.Range("D20").Value = "Farm to Table but also kinda Asian?"
$ws.Range("E20").Value = 43.673155291418801
$ws.Range("F20").Value = -79.3208615612443

# Update the current selection to reflect where the author left off editing
$ws.Range("D21").Select()
